$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.726.61'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '2.435.40'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.62'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.75'
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.169'
$ws.Range("E9").Value = '  +8.05%  '
$ws.Range("E10").Value = '  -2.16%  '
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("E12").Value = '  -5.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000177'
$ws.Range("E13").Value = '  +3.91%  '
$ws.Range("D14").Value = '68.624.07'
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = '2.882.83'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("E16").Value = '  -1.38%  '
$ws.Range("D17").Value = '2.437.03'
$ws.Range("E17").Value = '  -0.91%  '
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '339.10'
$ws.Range("E19").Value = '  +0.68%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.92'
$ws.Range("E22").Value = '  +2.26%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.70'
$ws.Range("E25").Value = '  +1.43%  '
$ws.Range("D26").Value = '2.564.76'
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.01'
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.20'
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.12'
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("E32").Value = '  +1.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '427.70'
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("E34").Value = '  -1.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.50'
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.97'
$ws.Range("E38").Value = '  +0.85%  '
$ws.Range("E39").Value = '  -2.67%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.50'
$ws.Range("E41").Value = '  +2.74%  '
$ws.Range("E42").Value = '  -1.96%  '
$ws.Range("E43").Value = '  -0.84%  '
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '130.94'
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("E46").Value = '  -1.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0720'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("E48").Value = '  -0.84%  '
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("E50").Value = '  +3.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0922'
$ws.Range("E51").Value = '  +0.61%  '
